# Update of league bases — swap mis-paired rows (id/home/away/odds columns
# were offset by one match), then drop the two trailing placeholder rows
# that had no real result yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B column (id) and F:AC columns (match data) need to be
# swapped with each other. Columns A, C, D, E (seq index / Div / Div name /
# Date) stay put on each row.
$rowPairs = @(
    @(20, 22),
    @(26, 27),
    @(35, 36),
    @(38, 39),
    @(43, 44),
    @(61, 62),
    @(73, 74),
    @(75, 76),
    @(102, 103),
    @(116, 117),
    @(118, 119),
    @(125, 126),
    @(127, 128),
    @(138, 139),
    @(156, 157),
    @(158, 159),
    @(168, 169),
    @(177, 178),
    @(183, 184),
    @(191, 192)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $bRange1 = $ws.Range("B$r1")
    $facRange1 = $ws.Range("F$r1" + ":AC$r1")
    $bRange2 = $ws.Range("B$r2")
    $facRange2 = $ws.Range("F$r2" + ":AC$r2")

    $bVal1 = $bRange1.Value()
    $facVal1 = $facRange1.Value()
    $bVal2 = $bRange2.Value()
    $facVal2 = $facRange2.Value()

    $bRange1.Value = $bVal2
    $facRange1.Value = $facVal2
    $bRange2.Value = $bVal1
    $facRange2.Value = $facVal1
}

# Remove the two trailing rows (194-195) that only had odds data and no
# final score / result yet — no longer part of the published dataset.
$ws.Rows("194:195").Delete()
